# Apply roster update to Sheet1: reorder/replace player rows 2-18,
# and remove the now-unused row 19 (table shrinks from 19 to 18 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for rows 2 through 18 (A: Player, B: Position, C: Team)
$data = @(
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Remove the last row (row 19), which is no longer part of the table
$ws.Rows.Item(19).Delete()
